$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 data
$ws.Range("A9").Value = "Corte-M4 Revision r0p1"
$ws.Range("B9").Value = "Technical Reference Manual"
$ws.Range("C9").Value = "Cortex-M4 Revision r0p1"
$ws.Range("D9").Value = "CortexM4_TRM_r0p1"

# Apply centered alignment style (style index 2) to B9 and C9 to match existing rows
$ws.Range("B9").HorizontalAlignment = -4108
$ws.Range("C9").HorizontalAlignment = -4108

# Update the view: scroll so column B is the left-most visible column, and select D9
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D9").Select()
